# Update cryptocurrency price/volume data and fix the Aave/EnergySwap row order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.603.62"
$ws.Range("E2").Value = "  +4.45%  "
$ws.Range("D3").Value = "2.659.10"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "568.19"
$ws.Range("E5").Value = "  +6.43%  "
$ws.Range("D6").Value = "146.52"
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").Value = "2.657.56"
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  +5.84%  "
$ws.Range("E12").Value = "  +7.09%  "
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").Value = "3.129.18"
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").Value = "60.559.49"
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").Value = "22.10"
$ws.Range("E16").Value = "  +6.95%  "
$ws.Range("E17").Value = "  +5.76%  "
$ws.Range("D18").Value = "2.657.40"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").Value = "343.42"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").Value = "10.46"
$ws.Range("E21").Value = "  +4.46%  "
$ws.Range("D22").Value = "6.40"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "66.28"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "7.41"
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("D29").Value = "0.0₃0807"
$ws.Range("E29").Value = "  +11.56%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "1.72"
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("E32").Value = "  +5.00%  "
$ws.Range("D33").Value = "159.21"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").Value = "19.20"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").Value = "4.12"
$ws.Range("E35").Value = "  +6.44%  "
$ws.Range("D36").Value = "0.897"
$ws.Range("E36").Value = "  +7.82%  "
$ws.Range("E37").Value = "  +6.18%  "
$ws.Range("D38").Value = "0.891"
$ws.Range("E38").Value = "  +8.92%  "
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("D40").Value = "37.45"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").Value = "298.39"
$ws.Range("E41").Value = "  +6.36%  "
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "0.0985"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "19.48"
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "128.44"
$ws.Range("E48").Value = "  +16.70%  "
$ws.Range("D49").Value = "10.71"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("E50").Value = "  +4.51%  "
$ws.Range("D51").Value = "18.78"
$ws.Range("E51").Value = "  +5.84%  "
